$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-OrganicFertilization")

$ws.Range("A1").Value = "Input [Sample Name]"
$ws.Range("AD1").Value = "Output [Sample Name]"
